# Reggio Calabria.xlsx update
# - Correct two historical "Nuovi casi" values (C414, C426) and let the
#   7-day rolling average formulas in column D recompute automatically.
# - Append daily rows for 30/07/2021 .. 19/08/2021 (rows 509-521) with new
#   case counts (column C) and their 7-day rolling averages (column D).
# - Append placeholder date-only rows for 12/08/2021 .. 31/08/2021
#   (rows 522-541) on every sheet (no data yet for those future days).
# - Update the active sheet / selections to reflect where the editor was
#   last working.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Data for the five sheets (sheetIndex -> daily "new count" column C
# values for rows 509..521, i.e. 2021-07-30 .. 2021-08-11)
# ---------------------------------------------------------------------
$newCaseValues = @(45,28,69,48,63,38,55,113,91,47,67,128,70)       # Nuovi casi
$deceasedValues = @(0,0,0,0,1,1,4,0,0,0,0,0,1)                     # Deceduti
$recoveredValues = @(0,72,2,6,11,37,24,17,20,19,14,14,0)           # Dimessi Guariti
$hospitalValues = @(17,23,24,26,26,27,27,28,27,30,33,33,37)        # Ricoveri
$icuValues = @(1,1,1,1,1,2,2,2,2,2,2,2,2)                          # Terapia

$sheetData = @{
    1 = $newCaseValues
    2 = $deceasedValues
    3 = $recoveredValues
    4 = $hospitalValues
    5 = $icuValues
}

# First new date serial is 44409 (2021-08-01) stored for row 511; row 509
# (44407) and row 510 (44408) already exist on every sheet but are still
# missing their C/D figures.
$firstDate = 44407
$firstRow = 509
$lastDataRow = 521      # last row that still receives a C/D figure
$lastBlankRow = 541     # sheets extend with date-only rows through here

for ($s = 1; $s -le 5; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $values = $sheetData[$s]

    # --- fill in C/D for the already-present rows 509 & 510, then create
    #     the new rows 511..521 with date + C + D -------------------------
    for ($r = $firstRow; $r -le $lastDataRow; $r++) {
        $dateSerial = $firstDate + ($r - $firstRow)
        $ws.Cells.Item($r, 1).Value = $dateSerial
        $ws.Cells.Item($r, 3).Value = $values[$r - $firstRow]
        $ws.Cells.Item($r, 4).Formula = "=AVERAGE(C" + ($r - 6) + ":C" + $r + ")"
    }

    # --- append the trailing date-only rows (no data yet) ----------------
    for ($r = ($lastDataRow + 1); $r -le $lastBlankRow; $r++) {
        $dateSerial = $firstDate + ($r - $firstRow)
        $ws.Cells.Item($r, 1).Value = $dateSerial
    }
}

# ---------------------------------------------------------------------
# Fix two historical data points on "Nuovi casi" - the 7-day rolling
# averages downstream recompute automatically.
# ---------------------------------------------------------------------
$wsNuoviCasi = $wb.Worksheets.Item(1)
$wsNuoviCasi.Cells.Item(414, 3).Value = 140
$wsNuoviCasi.Cells.Item(426, 3).Value = 109

# ---------------------------------------------------------------------
# View state: "Nuovi casi" becomes the active/selected sheet, with a
# fresh scroll position and selection on each sheet.
# ---------------------------------------------------------------------
$wsNuoviCasi.Activate() | Out-Null
$wsNuoviCasi.Range("C522").Select() | Out-Null

$win = $excel.ActiveWindow
$win.ScrollRow = 501
$win.ScrollColumn = 1

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("D521").Select() | Out-Null

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("D521").Select() | Out-Null

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("D521").Select() | Out-Null

$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("C521:D521").Select() | Out-Null

$wsNuoviCasi.Activate() | Out-Null
